# The three changed images live inside the document's headers/footers
# (the Pearson logo in the first-page & default footers, and the BTec
# logo in the first-page header). Word's InlineShape object model has
# no writable "Name" property (that's internal OOXML bookkeeping -
# wp:docPr/@name & pic:cNvPr/@name), so we round-trip the package
# through Document.WordOpenXML (exactly like Word's "paste as XML")
# and patch those two attributes on each picture in place.

$d = $word.ActiveDocument
$xml = $d.WordOpenXML

# 1) Footer "first page" (docPr id="3") - Pearson logo: image1.png -> image2.png
$old1 = '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="3" name="image1.png"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image1.png"/>'
$new1 = '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="3" name="image2.png"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image2.png"/>'
$xml = $xml.Replace($old1, $new1)

# 2) Footer "default" (docPr id="2") - Pearson logo: image1.png -> image2.png
$old2 = '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="2" name="image1.png"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image1.png"/>'
$new2 = '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="2" name="image2.png"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image2.png"/>'
$xml = $xml.Replace($old2, $new2)

# 3) Header "first page" (docPr id="1") - BTec logo: image2.jpg -> image1.jpg
$old3 = '<wp:docPr descr="BTec_Logo-Orange" id="1" name="image2.jpg"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="BTec_Logo-Orange" id="0" name="image2.jpg"/>'
$new3 = '<wp:docPr descr="BTec_Logo-Orange" id="1" name="image1.jpg"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="BTec_Logo-Orange" id="0" name="image1.jpg"/>'
$xml = $xml.Replace($old3, $new3)

$d.WordOpenXML = $xml
